$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values between row 2 and row 3 for the changed columns (D, M, N, O, P, R, S)

# Column D (Fecha) - stored as date serial numbers
$ws.Range("D2").Value2 = 44417
$ws.Range("D3").Value2 = 44235

# Column M (Volumen)
$ws.Range("M2").Value = 60
$ws.Range("M3").Value = 70

# Column N (Precio minimo)
$ws.Range("N2").Value = 26000
$ws.Range("N3").Value = 42000

# Column O (Precio maximo)
$ws.Range("O2").Value = 26000
$ws.Range("O3").Value = 42000

# Column P (Precio promedio ponderado)
$ws.Range("P2").Value = 26000
$ws.Range("P3").Value = 42000

# Column R (Origen)
$ws.Range("R2").Value = "Perú"
$ws.Range("R3").Value = "Región de Arica y Parinacota"

# Column S (Precio $/Kg)
$ws.Range("S2").Value = 1444
$ws.Range("S3").Value = 2333
